$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "61÷7=8, 5"
$t.Cell(1,2).Range.Text = "14÷8=1, 6"
$t.Cell(1,3).Range.Text = "38÷3=12, 2"
$t.Cell(1,4).Range.Text = "19÷2=9, 1"
$t.Cell(1,5).Range.Text = "93÷4=23, 1"

$t.Cell(5,1).Range.Text = "91÷6=15, 1"
$t.Cell(5,2).Range.Text = "98÷3=32, 2"
$t.Cell(5,3).Range.Text = "40÷7=5, 5"
$t.Cell(5,4).Range.Text = "43÷6=7, 1"
$t.Cell(5,5).Range.Text = "61÷4=15, 1"

$t.Cell(9,1).Range.Text = "64÷3=21, 1"
$t.Cell(9,2).Range.Text = "49÷6=8, 1"
$t.Cell(9,3).Range.Text = "49÷2=24, 1"
$t.Cell(9,4).Range.Text = "91÷5=18, 1"
$t.Cell(9,5).Range.Text = "28÷6=4, 4"

$t.Cell(13,1).Range.Text = "71÷2=35, 1"
$t.Cell(13,2).Range.Text = "15÷6=2, 3"
$t.Cell(13,3).Range.Text = "40÷4=10, 0"
$t.Cell(13,4).Range.Text = "97÷2=48, 1"
$t.Cell(13,5).Range.Text = "38÷4=9, 2"

$t.Cell(17,1).Range.Text = "18÷3=6, 0"
$t.Cell(17,2).Range.Text = "83÷3=27, 2"
$t.Cell(17,3).Range.Text = "76÷8=9, 4"
$t.Cell(17,4).Range.Text = "63÷6=10, 3"
$t.Cell(17,5).Range.Text = "22÷8=2, 6"

Write-Output "Done"
